$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "To stabilize the distribution we take logarithms of CO2.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To stabilize the distribution, we take logarithms of CO2.", 2
)
